$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Barrack Obama - ABSENT (no timestamp)
$ws.Range("A2").Value = "Barrack Obama"
$ws.Range("B2").Value = "ABSENT"
$ws.Range("C2").ClearContents()

# Row 3: Cyrus Gello Par - PRESENT - 16:46:57
$ws.Range("A3").Value = "Cyrus Gello Par"
$ws.Range("B3").Value = "PRESENT"
$ws.Range("C3").Value = "16:46:57"

# Row 4: Elon Musk - ABSENT (no timestamp)
$ws.Range("A4").Value = "Elon Musk"
$ws.Range("B4").Value = "ABSENT"
$ws.Range("C4").ClearContents()

# Row 5: Kiefer Tayawa - PRESENT - 16:46:45 (new row)
$ws.Range("A5").Value = "Kiefer Tayawa"
$ws.Range("B5").Value = "PRESENT"
$ws.Range("C5").Value = "16:46:45"

# Row 6: Roche Quejada - PRESENT - 16:46:50 (new row)
$ws.Range("A6").Value = "Roche Quejada"
$ws.Range("B6").Value = "PRESENT"
$ws.Range("C6").Value = "16:46:50"
